# Auto-generated: applies updated market-price data cell values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1188.1803
$ws.Range("I127").Value = 780
$ws.Range("J127").Value = 1216.8246
$ws.Range("K127").Value = 2340
$ws.Range("L127").Value = 3650.4738
$ws.Range("M127").Value = 2620
$ws.Range("N127").Value = -13570.4738
$ws.Range("H134").Value = 147700
$ws.Range("J134").Value = 147700
$ws.Range("L134").Value = 147700
$ws.Range("N134").Value = -157840
$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200
$ws.Range("H138").Value = 2187.2
$ws.Range("J138").Value = 2185.9656
$ws.Range("L138").Value = 6557.8968
$ws.Range("N138").Value = -16837.8968
$ws.Range("H139").Value = 62653.125
$ws.Range("J139").Value = 62653.125
$ws.Range("L139").Value = 62653.125
$ws.Range("N139").Value = -72933.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 38346
$ws.Range("J133").Value = 38346
$ws.Range("L133").Value = 38346
$ws.Range("N133").Value = -43406
$ws.Range("H134").Value = 27400
$ws.Range("J134").Value = 27400
$ws.Range("L134").Value = 27400
$ws.Range("N134").Value = -37540

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1273.625
$ws.Range("I94").Value = 979.82355
$ws.Range("J94").Value = 1987.1428
$ws.Range("K94").Value = 979.82355
$ws.Range("L94").Value = 1987.1428
$ws.Range("M94").Value = -528.82355
$ws.Range("N94").Value = -2889.1428
$ws.Range("H105").Value = 2100
$ws.Range("I105").Value = 2100
$ws.Range("K105").Value = 2100
$ws.Range("M105").Value = -353
$ws.Range("H107").Value = 3303.6667
$ws.Range("I107").Value = 3217.7334
$ws.Range("J107").Value = 3733.3333
$ws.Range("K107").Value = 3217.7334
$ws.Range("L107").Value = 3733.3333
$ws.Range("M107").Value = -1297.7334
$ws.Range("N107").Value = -7573.3333
$ws.Range("H135").Value = 68890
$ws.Range("J135").Value = 68890
$ws.Range("L135").Value = 68890
$ws.Range("N135").Value = -79030
$ws.Range("H140").Value = 52247.5
$ws.Range("J140").Value = 52247.5
$ws.Range("L140").Value = 52247.5
$ws.Range("N140").Value = -62607.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19233326
$ws.Range("I31").Value = 28573178
$ws.Range("K31").Value = 28573178
$ws.Range("M31").Value = -28572883
$ws.Range("H34").Value = 19233326
$ws.Range("I34").Value = 28573178
$ws.Range("K34").Value = 28573178
$ws.Range("M34").Value = -28572976
$ws.Range("H138").Value = 40979
$ws.Range("J138").Value = 40979
$ws.Range("L138").Value = 40979
$ws.Range("N138").Value = -51259
$ws.Range("H140").Value = 81876.16
$ws.Range("J140").Value = 81876.16
$ws.Range("L140").Value = 81876.16
$ws.Range("N140").Value = -92236.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 21333.334
$ws.Range("I70").Value = 30000
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 90000
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -89685
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 21333.334
$ws.Range("I73").Value = 30000
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 90000
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -88908
$ws.Range("N73").Value = -14184
$ws.Range("H76").Value = 2950
$ws.Range("I76").Value = 1900
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 5700
$ws.Range("L76").Value = 12000
$ws.Range("M76").Value = -5317
$ws.Range("N76").Value = -12766
$ws.Range("H79").Value = 2950
$ws.Range("I79").Value = 1900
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 5700
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = -4374
$ws.Range("N79").Value = -14652
$ws.Range("H82").Value = 13502.5
$ws.Range("J82").Value = 14717.143
$ws.Range("L82").Value = 44151.429
$ws.Range("N82").Value = -44963.429
$ws.Range("H85").Value = 13502.5
$ws.Range("J85").Value = 14717.143
$ws.Range("L85").Value = 44151.429
$ws.Range("N85").Value = -46959.429
$ws.Range("H88").Value = 104000
$ws.Range("J88").Value = 104000
$ws.Range("L88").Value = 312000
$ws.Range("N88").Value = -312856
$ws.Range("H91").Value = 104000
$ws.Range("J91").Value = 104000
$ws.Range("L91").Value = 312000
$ws.Range("N91").Value = -314964
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 3822
$ws.Range("J94").Value = 3975
$ws.Range("L94").Value = 11925
$ws.Range("N94").Value = -13277
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H96").Value = 14600
$ws.Range("J96").Value = 5750
$ws.Range("L96").Value = 17250
$ws.Range("N96").Value = -21368
$ws.Range("H97").Value = 500
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H99").Value = 5481.25
$ws.Range("J99").Value = 10000
$ws.Range("L99").Value = 30000
$ws.Range("N99").Value = -34492
$ws.Range("H100").Value = 21500
$ws.Range("I100").Value = 50000
$ws.Range("K100").Value = 150000
$ws.Range("M100").Value = -149189
$ws.Range("H101").Value = 5565.489
$ws.Range("J101").Value = 5565.489
$ws.Range("L101").Value = 16696.467
$ws.Range("N101").Value = -21564.467
$ws.Range("H103").Value = 2145.6667
$ws.Range("J103").Value = 5500
$ws.Range("L103").Value = 16500
$ws.Range("N103").Value = -18258
$ws.Range("H107").Value = 660.1667
$ws.Range("I107").Value = 1172.5
$ws.Range("J107").Value = 404
$ws.Range("K107").Value = 3517.5
$ws.Range("L107").Value = 1212
$ws.Range("M107").Value = -1597.5
$ws.Range("N107").Value = -5052
$ws.Range("H137").Value = 55557536
$ws.Range("J137").Value = 111114130
$ws.Range("L137").Value = 333342390
$ws.Range("N137").Value = -333352590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2707.2144
$ws.Range("I126").Value = 2151
$ws.Range("K126").Value = 6453
$ws.Range("M126").Value = -3983
$ws.Range("H132").Value = 2958.75
$ws.Range("I132").Value = 1783.1818
$ws.Range("J132").Value = 4395.5557
$ws.Range("K132").Value = 5349.5454
$ws.Range("L132").Value = 13186.6671
$ws.Range("M132").Value = -2819.5454
$ws.Range("N132").Value = -18246.6671
$ws.Range("H135").Value = 47151
$ws.Range("J135").Value = 47151
$ws.Range("L135").Value = 47151
$ws.Range("N135").Value = -57291
$ws.Range("H140").Value = 34001.25
$ws.Range("J140").Value = 34001.25
$ws.Range("L140").Value = 34001.25
$ws.Range("N140").Value = -44361.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 78919.664
$ws.Range("J134").Value = 78919.664
$ws.Range("L134").Value = 78919.664
$ws.Range("N134").Value = -89059.664
$ws.Range("H137").Value = 83964.5
$ws.Range("J137").Value = 83964.5
$ws.Range("L137").Value = 83964.5
$ws.Range("N137").Value = -94164.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4157.8335
$ws.Range("I122").Value = 2716.8333
$ws.Range("J122").Value = 5598.8335
$ws.Range("K122").Value = 8150.499899999999
$ws.Range("L122").Value = 16796.5005
$ws.Range("M122").Value = -5700.499899999999
$ws.Range("N122").Value = -21696.5005
$ws.Range("H133").Value = 29846.666
$ws.Range("J133").Value = 29846.666
$ws.Range("L133").Value = 29846.666
$ws.Range("N133").Value = -39966.666
$ws.Range("H135").Value = 44009.453
$ws.Range("J135").Value = 44009.453
$ws.Range("L135").Value = 44009.453
$ws.Range("N135").Value = -54149.453
$ws.Range("H139").Value = 57280
$ws.Range("J139").Value = 57280
$ws.Range("L139").Value = 57280
$ws.Range("N139").Value = -67560
$ws.Range("H141").Value = 69499.84
$ws.Range("J141").Value = 69499.84
$ws.Range("L141").Value = 69499.84
$ws.Range("N141").Value = -79859.84
